# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values to the Siren_Profits workbook
# (columns H..N: currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ), LeveProfit(NQ/HQ))
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ALC_updates = [ordered]@{
    "H29" = 3366.3333
    "I29" = 99
    "K29" = 297
    "M29" = -16
    "H38" = 31.8
    "I38" = 33.666668
    "J38" = 15
    "K38" = 101.000004
    "L38" = 45
    "M38" = 270.999996
    "N38" = -789
    "H41" = 574.7619
    "I41" = 342
    "J41" = 786.36365
    "K41" = 342
    "L41" = 786.36365
    "M41" = 98
    "N41" = -1666.36365
    "H58" = 2667.6667
    "I58" = 203
    "J58" = 3900
    "K58" = 609
    "L58" = 11700
    "M58" = -459
    "N58" = -12000
    "H125" = 4596.643
    "J125" = 2683
    "L125" = 24147
    "N125" = -29067
    "H138" = 151658.6
    "J138" = 4512.8438
    "L138" = 13538.5314
    "N138" = -23818.5314
}
foreach ($ref in $ALC_updates.Keys) {
    $ws.Range($ref).Value = $ALC_updates[$ref]
}

$ws = $wb.Worksheets.Item("ARM")
$ARM_updates = [ordered]@{
    "H2" = 3347.3235
    "I2" = 3188.88
    "K2" = 3188.88
    "M2" = -3075.88
    "H32" = 608786.7
    "I32" = 608786.7
    "K32" = 608786.7
    "M32" = -608499.7
    "H45" = 5905
    "I45" = 5125.625
    "J45" = 8399
    "K45" = 5125.625
    "L45" = 8399
    "M45" = -4748.625
    "N45" = -9153
    "H116" = 3347.3235
    "I116" = 3188.88
    "K116" = 3188.88
    "M116" = -894.8800000000001
    "H132" = 2087.22
    "I132" = 1143.8158
    "J132" = 5074.6665
    "K132" = 3431.4474
    "L132" = 15223.9995
    "M132" = -901.4474
    "N132" = -20283.9995
    "H138" = 72469
    "J138" = 72469
    "L138" = 72469
    "N138" = -82749
}
foreach ($ref in $ARM_updates.Keys) {
    $ws.Range($ref).Value = $ARM_updates[$ref]
}

$ws = $wb.Worksheets.Item("BSM")
$BSM_updates = [ordered]@{
    "H3" = 3347.3235
    "I3" = 3188.88
    "K3" = 3188.88
    "M3" = -3074.88
    "H74" = 48549.668
    "J74" = 48549.668
    "L74" = 48549.668
    "N74" = -50421.668
    "H77" = 48549.668
    "J77" = 48549.668
    "L77" = 145649.004
    "N77" = -155009.004
    "H105" = 14171.833
    "I105" = 20458.5
    "K105" = 20458.5
    "M105" = -18711.5
    "H130" = 89999
    "J130" = 89999
    "L130" = 89999
    "N130" = -100039
    "H134" = 1810.0714
    "I134" = 1374.3158
    "K134" = 4122.9474
    "M134" = -1587.9474
    "H140" = 116155.8
    "J140" = 120194.75
    "L140" = 120194.75
    "N140" = -130554.75
}
foreach ($ref in $BSM_updates.Keys) {
    $ws.Range($ref).Value = $BSM_updates[$ref]
}

$ws = $wb.Worksheets.Item("CRP")
$CRP_updates = [ordered]@{
    "H31" = 2538.7659
    "I31" = 1873.5135
    "J31" = 5000.2
    "K31" = 1873.5135
    "L31" = 5000.2
    "M31" = -1578.5135
    "N31" = -5590.2
    "H34" = 2538.7659
    "I34" = 1873.5135
    "J34" = 5000.2
    "K34" = 1873.5135
    "L34" = 5000.2
    "M34" = -1671.5135
    "N34" = -5404.2
    "H58" = 2762.913
    "I58" = 2542.6875
    "J58" = 3266.2856
    "K58" = 2542.6875
    "L58" = 3266.2856
    "M58" = -2339.6875
    "N58" = -3672.2856
    "H105" = 7253.5
    "I105" = 8014.579
    "J105" = 2433.3333
    "K105" = 8014.579
    "L105" = 2433.3333
    "M105" = -6267.579
    "N105" = -5927.3333
    "H130" = 60000
    "J130" = 60000
    "L130" = 60000
    "N130" = -70040
    "H136" = 2762.913
    "I136" = 2542.6875
    "J136" = 3266.2856
    "K136" = 7628.0625
    "L136" = 9798.856800000001
    "M136" = -5078.0625
    "N136" = -14898.8568
    "H141" = 341341.44
    "J141" = 415540.22
    "L141" = 415540.22
    "N141" = -425900.22
}
foreach ($ref in $CRP_updates.Keys) {
    $ws.Range($ref).Value = $CRP_updates[$ref]
}

$ws = $wb.Worksheets.Item("CUL")
$CUL_updates = [ordered]@{
    "H55" = 10794.366
    "J55" = 11237.667
    "L55" = 33713.001
    "N55" = -34067.001
    "H131" = 34486970
    "I131" = 83340710
    "J131" = 1970.6471
    "K131" = 250022130
    "L131" = 5911.9413
    "M131" = -250017090
    "N131" = -15991.9413
    "H139" = 5717857
    "I139" = 20000000
    "J139" = 5000
    "K139" = 60000000
    "L139" = 15000
    "M139" = -59994860
    "N139" = -25280
    "H140" = 1252207
    "I140" = 1430736.6
    "K140" = 4292209.800000001
    "M140" = -4287029.800000001
}
foreach ($ref in $CUL_updates.Keys) {
    $ws.Range($ref).Value = $CUL_updates[$ref]
}

$ws = $wb.Worksheets.Item("GSM")
$GSM_updates = [ordered]@{
    "H102" = 7831.48
    "I102" = 8466.048000000001
    "K102" = 8466.048000000001
    "M102" = -6844.048000000001
    "H131" = 27699.5
    "J131" = 27699.5
    "L131" = 27699.5
    "N131" = -37779.5
}
foreach ($ref in $GSM_updates.Keys) {
    $ws.Range($ref).Value = $GSM_updates[$ref]
}

$ws = $wb.Worksheets.Item("LTW")
$LTW_updates = [ordered]@{
    "H16" = 3966
    "I16" = 3808.0908
    "J16" = 4313.4
    "K16" = 3808.0908
    "L16" = 4313.4
    "M16" = -3638.0908
    "N16" = -4653.4
    "H22" = 738.625
    "I22" = 707.3333
    "K22" = 707.3333
    "M22" = -412.3333
    "H27" = 738.625
    "I27" = 707.3333
    "K27" = 707.3333
    "M27" = -600.3333
    "H46" = 2214.55
    "J46" = 4399.5
    "L46" = 4399.5
    "N46" = -4775.5
    "H61" = 21937.73
    "I61" = 1883.6923
    "K61" = 1883.6923
    "M61" = -1681.6923
    "H113" = 21937.73
    "I113" = 1883.6923
    "K113" = 1883.6923
    "M113" = 286.3077000000001
    "H122" = 7491.9473
    "I122" = 7734.7
    "J122" = 7222.222
    "K122" = 23204.1
    "L122" = 21666.666
    "M122" = -20754.1
    "N122" = -26566.666
}
foreach ($ref in $LTW_updates.Keys) {
    $ws.Range($ref).Value = $LTW_updates[$ref]
}

$ws = $wb.Worksheets.Item("WVR")
$WVR_updates = [ordered]@{
    "H74" = 113625.5
    "I74" = 5000
    "J74" = 149834
    "K74" = 5000
    "L74" = 149834
    "M74" = -4064
    "N74" = -151706
    "H77" = 113625.5
    "I77" = 5000
    "J77" = 149834
    "K77" = 15000
    "L77" = 449502
    "M77" = -10320
    "N77" = -458862
    "H81" = 11642.091
    "I81" = 26812.75
    "J81" = 2973.1428
    "K81" = 53625.5
    "L81" = 5946.2856
    "M81" = -52564.5
    "N81" = -8068.2856
    "H84" = 11642.091
    "I84" = 26812.75
    "J84" = 2973.1428
    "K84" = 268127.5
    "L84" = 29731.428
    "M84" = -262823.5
    "N84" = -40339.428
}
foreach ($ref in $WVR_updates.Keys) {
    $ws.Range($ref).Value = $WVR_updates[$ref]
}
